# EnterpriseJavaTimeLog.xlsx edit
#
# Commit: "Adds capability to search using multiple criteria (dao) and to
# display a story on the profile"
#
# Concretely this:
#   1. Inserts two new blank rows into the time-log table (after the
#      existing blank-row block, before the "Issues/Loose Ends" notes
#      section), shifting all later rows down by two.
#   2. Fills in the first of the newly-available blank rows (row 41) with
#      a new time-log entry: date 2019-03-12, 2 hours, and a description
#      of the work done.
#   3. Updates the active selection to reflect where the user ended up
#      (D42) after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows just below row 43 (i.e. before the old row 44),
# pushing the "Issues/Loose Ends" section and everything after it down by
# two rows.
$ws.Rows("44:45").Insert()

# Populate the newly available row 41 with the new time-log entry.
$ws.Range("A41").Value = 43536
$ws.Range("B41").Value = 2
$ws.Range("D41").Value = "Indie Project: Trying Paula's example function for retrieving an entity based on its own characteristics and that of another entity. Tested it in StoryDaoTest. Used it for the profile display servlet/jsp."

# The wrapped description makes the row taller, same as the other
# multi-line entries above it.
$ws.Rows(41).RowHeight = 30

# Reflect the user's resulting cursor position.
$ws.Range("D42").Select()
